# Insert a new weekly data row above row 133, shifting existing rows down.
# The new row duplicates the surrounding metadata (market, region, category,
# quality, unit, origin, classification) and introduces a fresh date plus
# updated volume / price figures, exactly like every other week's entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(133).Insert()

$ws.Cells.Item(133, 1).Value  = 8
$ws.Cells.Item(133, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(133, 3).Value  = "Coquimbo"
$ws.Cells.Item(133, 4).Value  = 44960
$ws.Cells.Item(133, 5).Value  = 4
$ws.Cells.Item(133, 6).Value  = 100112001
$ws.Cells.Item(133, 7).Value  = "Berenjena"
$ws.Cells.Item(133, 8).Value  = "Sin especificar"
$ws.Cells.Item(133, 9).Value  = "Primera"
$ws.Cells.Item(133, 10).Value = 520
$ws.Cells.Item(133, 11).Value = 11800
$ws.Cells.Item(133, 12).Value = 12000
$ws.Cells.Item(133, 13).Value = 11900
$ws.Cells.Item(133, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(133, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(133, 16).Value = 298
$ws.Cells.Item(133, 17).Value = 40
$ws.Cells.Item(133, 18).Value = "Hortaliza"
